# Add a second "video link" textbox to slide 1, right below the existing
# YouTube-link textbox ("TextBox 3"), mirroring its exact formatting
# (no fill, wrap="none", spAutoFit, etc.) and set its own text/position.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the existing "TextBox 3" shape (the one holding the first
# YouTube link) so the new shape inherits identical shape/text formatting.
$source = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 3") {
        $source = $candidate
    }
}

# Duplicate it -- this clones all shape/body formatting (noFill, wrap,
# spAutoFit, lstStyle, run formatting, ...) exactly, so we only need to
# change what actually differs for the new box: its name, position/size,
# and text.
$newShapes = $source.Duplicate()
$newBox = $newShapes.Item(1)

$newBox.Name = "TextBox 4"

# Position/size target values are expressed in EMU in the OOXML; PowerPoint's
# Shape.Left/Top/Width/Height are in points, so convert (1 pt = 12700 EMU).
$newBox.Left = 5155512 / 12700
$newBox.Top = 4579167 / 12700
$newBox.Width = 3384260 / 12700
$newBox.Height = 369332 / 12700

$newBox.TextFrame.TextRange.Text = "https://youtu.be/nk3bkaz5fKg"
